# This script reproduces the edits described in the diff:
#  - Rows 10, 11, 12 on sheet "bets" get their bet-detail columns
#    (C, E, H, I, J, K, L) rewritten; the D/F/M formula cells
#    recalculate automatically.
#  - A brand-new, otherwise-empty row 18 is added, carrying the
#    date-column format in B18 and the percentage-column format in M18.
#  - The active selection ends on E8, matching the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 10 ----------------------------------------------------------
$ws.Range("C10").Value = 1
$ws.Range("E10").Value = 1140
$ws.Range("H10").Value = "HALO WORLDS 2023"
$ws.Range("I10").Value = "OPTIC"
$ws.Range("J10").Value = "GANA SERIE"
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0

# ---- Row 11 ----------------------------------------------------------
$ws.Range("E11").Value = 690
$ws.Range("H11").Value = "VALORANT CHAMPIONS LATAM"
$ws.Range("I11").Value = "KRU"

# ---- Row 12 ----------------------------------------------------------
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = -5000
$ws.Range("H12").Value = "WORLDS 2023"
$ws.Range("I12").Value = "LOUD"
$ws.Range("J12").Value = "GANA 1 MAPA EN LA SERIE"
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1

# ---- New row 18 --------------------------------------------------------
# Copy the existing date / percentage formatting down onto the new row
# so the blank cells keep the same look as the rest of the table.
$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M18").Style = "Porcentaje"

$excel.CutCopyMode = $false

# ---- Selection ---------------------------------------------------------
$ws.Range("E8").Select()
